# Generate Report for Handback
#
# Updates the "Latest Handback DateTime" (column K) for the
# 3f2c5e87-344b-4b17-a884-2bd3de23af6e file row (row 2) on both the
# "zh-cn" and "de-de" localization-status sheets, reflecting a newer
# handback timestamp recorded by the report generator.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2 on each sheet corresponds to 3f2c5e87-344b-4b17-a884-2bd3de23af6e;
# column K is "Latest Handback DateTime".
$wsZhCn.Range("K2").Value = "2016-10-17 16:49:25"
$wsDeDe.Range("K2").Value = "2016-10-17 16:50:02"
